# Add the "2022-Q4" sheet (placed right after "总计") with fund holding data,
# and update the "总计" (summary) sheet with the new quarter's row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q4" worksheet by copying an existing quarter sheet
#    (this keeps column widths / styles / borders identical to its siblings),
#    then trim it down to the 3 fund rows we actually need and overwrite the
#    values.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2022-Q3")

$template.Copy($null, $total)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The template sheet has 8 data rows (rows 2-9); we only need 3 (rows 2-4).
$newSheet.Rows("5:9").Delete()

# Columns B:G hold text values (fund code / name / percentages formatted as
# strings) -- force text storage so things like "010922" or "15.00" keep
# their exact formatting, then clear the temporary format so the cells end
# up with no explicit style (matching the other data cells on the sheet).
$dataText = $newSheet.Range("B2:G4")
$dataText.NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "010922"
$newSheet.Range("C2").Value = "朱雀匠心一年持有期混合"
$newSheet.Range("D2").Value = "15.00"
$newSheet.Range("E2").Value = "93.80"
$newSheet.Range("F2").Value = "2.84"
$newSheet.Range("G2").Value = "0.4260"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "004266"
$newSheet.Range("C3").Value = "招商沪港深科技创新主题精选灵活配置混合A"
$newSheet.Range("D3").Value = "0.93"
$newSheet.Range("E3").Value = "90.63"
$newSheet.Range("F3").Value = "2.09"
$newSheet.Range("G3").Value = "0.0194"
$newSheet.Range("H3").Value = 10

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "010754"
$newSheet.Range("C4").Value = "招商沪港深科技创新主题精选灵活配置混合C"
$newSheet.Range("D4").Value = "0.25"
$newSheet.Range("E4").Value = "90.63"
$newSheet.Range("F4").Value = "2.09"
$newSheet.Range("G4").Value = "0.0052"
$newSheet.Range("H4").Value = 10

$dataText.ClearFormats()

# ---------------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a new row for 2022-Q4 above the existing
#    2022-Q3 row, shifting the rest of the table down.
# ---------------------------------------------------------------------------
$total.Rows("2:2").Insert()

# Restore column A's style (the inserted row doesn't inherit it automatically)
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Clear the stray formatting the insert copied into B2:D2 so the row matches
# the plain (style-less) look of the other data rows.
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.45

# The shifted rows keep their old index counter (0,1,2,3); renumber column A
# back into a continuous 0..4 sequence.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4

# Keep the focus on the summary sheet, same as before the edit.
$total.Activate()

Write-Host "2022-Q4 sheet added and 总计 summary updated"
